# Generate Report for Handback
#
# The CI run that produced this workbook has now handed back the file
# "9ef8f679-eb76-4707-8fa4-8517cddb2476.md" for both locales (zh-cn and
# de-de): it was previously "Ready for handoff" and is now
# "Handed back: in sync with en-US", with a fresh handback timestamp
# recorded in the "Latest Handback DateTime" column.
#
# Update the three worksheets accordingly:
#   - Overview           : row for the file -> zh-cn / de-de status cells
#   - zh-cn (per-locale)  : Status cell + Latest Handback DateTime cell
#   - de-de (per-locale)  : Status cell + Latest Handback DateTime cell

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusHandedBack
$wsZhCn.Range("G3").Value = "2016-03-10 22:59:27"

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusHandedBack
$wsDeDe.Range("G3").Value = "2016-03-10 22:59:43"
